# Add the new "ODI Batting Extra" worksheet (4th sheet) and populate it
# with the PlayerPerformance batting-extra data.

$wb = $excel.ActiveWorkbook

# --- Create the sheet, name it, and move it to the end of the tab order ---
$ws = $wb.Worksheets.Add()
$ws.Name = "ODI Batting Extra"
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch a live reference to the sheet by name now that it has moved -
# the old $ws handle can end up pointing at the wrong sheet after Move().
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# --- Header row: reuse the bold/centered/bordered header formatting that
# the other sheets in this workbook already use, then fill in the labels. ---
$headerSource = $wb.Worksheets.Item("ODI Bowling").Range("A1:F1")
$headerSource.Copy($null)
$ws.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data rows ---
# MATCH_CODE (A), NUM_4 (C), NUM_6 (D) and PERCENT_RUNS_OF_TOTAL (E) are stored
# as text, so force a text number-format before assigning numeric-looking
# strings (otherwise Excel auto-converts them to numbers/percentages).
# BATTING_POSITION (B) is a genuine number. MAN_OF_MATCH (F) is plain text.

$data = @(
    @("4108", "",   "",  "",  "",      "NO"),
    @("4115", 10,   "0", "1", "2.96%", "NO"),
    @("4123", 9,    "0", "0", "1.52%", "NO"),
    @("4125", 9,    "1", "0", "3.24%", "NO"),
    @("4166", 9,    "2", "1", "8.88%", "NO"),
    @("4167", "",   "",  "",  "",      "NO"),
    @("4168", 10,   "0", "0", "2.09%", "NO")
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $r + 2
    $values = $data[$r]

    $a = $ws.Cells.Item($row, 1)
    $a.NumberFormat = "@"
    $a.Value = $values[0]

    $ws.Cells.Item($row, 2).Value = $values[1]

    $c = $ws.Cells.Item($row, 3)
    $c.NumberFormat = "@"
    $c.Value = $values[2]

    $d = $ws.Cells.Item($row, 4)
    $d.NumberFormat = "@"
    $d.Value = $values[3]

    $e = $ws.Cells.Item($row, 5)
    $e.NumberFormat = "@"
    $e.Value = $values[4]

    $ws.Cells.Item($row, 6).Value = $values[5]
}
